$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44547
$ws.Range("J2").Value = 200

# Row 3
$ws.Range("D3").Value = 44568
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 15000
$ws.Range("L3").Value = 16000
$ws.Range("M3").Value = 15500
$ws.Range("P3").Value = 861

# Row 4
$ws.Range("D4").Value = 44557
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 13000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 13500
$ws.Range("P4").Value = 750
